# Sourced male and female audio inputs, selector switch and power button
# Adds a new "Audio Input" BOM worksheet after "Power", populated with the
# new parts, and moves the "active sheet" selection from Power to the new
# Audio Input sheet (matching how Excel leaves the previously-active sheet's
# selection parked on its header row).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Park the selection on the (soon to be previously active) Power sheet
#    header row, the way Excel leaves a sheet you navigate away from.
# ---------------------------------------------------------------------
$powerSheet = $wb.Worksheets.Item("Power")
$powerSheet.Activate()
$powerSheet.Range("A1:F1").Select()

# ---------------------------------------------------------------------
# 2. Add the new "Audio Input" sheet right after "Power" (the last sheet).
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "Audio Input"

# ---------------------------------------------------------------------
# 3. Copy the header-row formatting (bold + bottom border) from Power,
#    which is identical to what the new sheet needs.
# ---------------------------------------------------------------------
$powerSheet.Range("A1:F1").Copy()
$newSheet.Range("A1:F1").PasteSpecial(-4122) | Out-Null
$newSheet.Application.CutCopyMode = $false

# ---------------------------------------------------------------------
# 4. Header row text (reuses existing shared strings from other sheets).
# ---------------------------------------------------------------------
$newSheet.Range("A1").Value = "Component"
$newSheet.Range("B1").Value = "Value /Purpose"
$newSheet.Range("C1").Value = "Price "
$newSheet.Range("D1").Value = "Link "
$newSheet.Range("E1").Value = "Case Code (mm)"
$newSheet.Range("F1").Value = "Datasheet"

# ---------------------------------------------------------------------
# 5. Data rows: male jack, female jack, selector switch, on/off switch.
# ---------------------------------------------------------------------
$newSheet.Range("A2").Value = "SP-3533-02"
$newSheet.Range("B2").Value = "Male Jack"
$newSheet.Range("C2").Value = 2.98
$newSheet.Range("D2").Value = "https://www.mouser.ca/ProductDetail/CUI-Devices/SP-3533-02?qs=%252BEew9%252B0nqrCkfyt%2FFhI%252B5A%3D%3D"
$newSheet.Range("F2").Value = "SP-3533-02 Datasheet (PDF)"

$newSheet.Range("A3").Value = "STX-3000"
$newSheet.Range("B3").Value = "Female Jack"
$newSheet.Range("C3").Value = 1.02
$newSheet.Range("D3").Value = "https://www.mouser.ca/ProductDetail/Kycon/STX-3000?qs=kjZ2mQLP346Nbz1X9BOzfg%3D%3D"
$newSheet.Range("F3").Value = "STX-3000 Datasheet (PDF)"

$newSheet.Range("A4").Value = "OS102011MS2QN1"
$newSheet.Range("B4").Value = "Slide Selector Switch"
$newSheet.Range("C4").Value = 0.538
$newSheet.Range("D4").Value = "https://www.mouser.ca/ProductDetail/CK/OS102011MS2QN1?qs=WtljUlYws5RvQ1hEv876nQ%3D%3D"
$newSheet.Range("F4").Value = "OS102011MS2QN1 Datasheet (PDF)"

$newSheet.Range("A5").Value = "RA11131123"
$newSheet.Range("B5").Value = "On/OFF switch"
$newSheet.Range("C5").Value = 0.925
$newSheet.Range("D5").Value = "https://www.mouser.ca/ProductDetail/E-Switch/RA11131123?qs=QtyuwXswaQhc6OwdGDJDiQ%3D%3D"
$newSheet.Range("F5").Value = "RA11131123 Datasheet (PDF)"

# ---------------------------------------------------------------------
# 6. Datasheet hyperlinks on column F (display text mirrors the PDF URL).
# ---------------------------------------------------------------------
$newSheet.Hyperlinks.Add($newSheet.Range("F2"), "https://www.mouser.ca/datasheet/2/670/sp_3533_02-1779105.pdf", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "https://www.mouser.ca/datasheet/2/670/sp_3533_02-1779105.pdf") | Out-Null
$newSheet.Hyperlinks.Add($newSheet.Range("F3"), "https://www.mouser.ca/datasheet/2/222/STX3000-334650.pdf", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "https://www.mouser.ca/datasheet/2/222/STX3000-334650.pdf") | Out-Null
$newSheet.Hyperlinks.Add($newSheet.Range("F4"), "https://www.mouser.ca/datasheet/2/60/os-1841429.pdf", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "https://www.mouser.ca/datasheet/2/60/os-1841429.pdf") | Out-Null
$newSheet.Hyperlinks.Add($newSheet.Range("F5"), "https://www.mouser.ca/datasheet/2/140/ESCH_S_A0005379088_1-2548267.pdf", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "https://www.mouser.ca/datasheet/2/140/ESCH_S_A0005379088_1-2548267.pdf") | Out-Null

# Re-apply the Hyperlink cell style so it matches the rest of the workbook
# (Hyperlinks.Add already colors/underlines the cell; this just pins the
# canonical named style rather than an ad-hoc copy of it).
$newSheet.Range("F2:F5").Style = "Hyperlink"

# ---------------------------------------------------------------------
# 7. Column widths roughly matching the authored sheet (bestfit-like).
# ---------------------------------------------------------------------
$newSheet.Columns.Item(1).ColumnWidth = 16
$newSheet.Columns.Item(2).ColumnWidth = 13.6

# ---------------------------------------------------------------------
# 8. Page setup + final active sheet/selection state.
# ---------------------------------------------------------------------
$newSheet.PageSetup.Orientation = 1

$newSheet.Activate()
$newSheet.Range("F12").Select()
